# Add 2022-Q3 data: insert a new "2022-Q3" sheet (fund holdings) between
# "总计" and "2022-Q2", and add a corresponding summary row in "总计".

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) "总计" sheet: insert a new row for 2022-Q3 above the existing data,
#    pushing the 2022-Q2 / 2022-Q1 rows down by one.
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("总计")

# Push row 3 -> row 4, copying formatting (so styles match exactly),
# then row 2 -> row 3, then overwrite row 2 with the new 2022-Q3 figures.
# (Restrict the copy to the used A:D columns -- copying whole Rows would
# stamp style refs across all 16384 columns.)
$summary.Range("A3:D3").Copy()
$summary.Range("A4:D4").PasteSpecial(-4122)   # xlPasteFormats

$summary.Range("A2:D2").Copy()
$summary.Range("A3:D3").PasteSpecial(-4122)   # xlPasteFormats

$summary.Cells.Item(4,1).Value2 = 2
$summary.Cells.Item(4,2).Value2 = "2022-Q1"
$summary.Cells.Item(4,3).Value2 = 3
$summary.Cells.Item(4,4).Value2 = 0.7

$summary.Cells.Item(3,1).Value2 = 1
$summary.Cells.Item(3,2).Value2 = "2022-Q2"
$summary.Cells.Item(3,3).Value2 = 5
$summary.Cells.Item(3,4).Value2 = 0.45

$summary.Cells.Item(2,1).Value2 = 0
$summary.Cells.Item(2,2).Value2 = "2022-Q3"
$summary.Cells.Item(2,3).Value2 = 5
$summary.Cells.Item(2,4).Value2 = 0.09

# ---------------------------------------------------------------------
# 2) Insert the new "2022-Q3" sheet right before "2022-Q2" by duplicating
#    the "2022-Q2" sheet (keeps sheetPr/margins/styles identical) and
#    then overwriting its data with the new quarter's figures.
# ---------------------------------------------------------------------
$q2 = $wb.Worksheets.Item("2022-Q2")
$q2.Copy($q2)

# Re-fetch sheet references: inserting a sheet invalidates previously
# held handles (reading from them afterwards silently yields blanks).
$q3 = $wb.Worksheets.Item("2022-Q2 (2)")
$q3.Name = "2022-Q3"

# Header row (plain text, not numeric-looking -> safe to set directly).
$q3.Cells.Item(1,2).Value2 = "基金代码"
$q3.Cells.Item(1,3).Value2 = "基金名称"
$q3.Cells.Item(1,4).Value2 = "基金规模"
$q3.Cells.Item(1,5).Value2 = "股票总仓位"
$q3.Cells.Item(1,6).Value2 = "仓位占比"
$q3.Cells.Item(1,7).Value2 = "持有市值(亿元)"
$q3.Cells.Item(1,8).Value2 = "仓位排名"

# Numeric columns (A = row index, H = rank) -> plain values.
$q3.Cells.Item(2,1).Value2 = 0
$q3.Cells.Item(3,1).Value2 = 1
$q3.Cells.Item(4,1).Value2 = 2
$q3.Cells.Item(5,1).Value2 = 3
$q3.Cells.Item(6,1).Value2 = 4

$q3.Cells.Item(2,8).Value2 = 8
$q3.Cells.Item(3,8).Value2 = 2
$q3.Cells.Item(4,8).Value2 = 3
$q3.Cells.Item(5,8).Value2 = 8
$q3.Cells.Item(6,8).Value2 = 6

# Text columns B-G hold numeric-looking strings (fund codes, percentages,
# etc.) that must stay text (t="inlineStr"/shared string), not be
# silently converted to numbers by Excel. Route them through a text
# formula, then collapse the formula to its literal value, which avoids
# both the numeric auto-conversion AND the "number stored as text"
# quote-prefix styling that a direct text assignment would pick up.
$q3.Cells.Item(2,2).Formula = '="012245"'
$q3.Cells.Item(2,3).Formula = '="广发金融地产精选股票C"'
$q3.Cells.Item(2,4).Formula = '="0.55"'
$q3.Cells.Item(2,5).Formula = '="85.52"'
$q3.Cells.Item(2,6).Formula = '="4.62"'
$q3.Cells.Item(2,7).Formula = '="0.0254"'

$q3.Cells.Item(3,2).Formula = '="009999"'
$q3.Cells.Item(3,3).Formula = '="东方中国红利混合"'
$q3.Cells.Item(3,4).Formula = '="0.51"'
$q3.Cells.Item(3,5).Formula = '="79.37"'
$q3.Cells.Item(3,6).Formula = '="4.73"'
$q3.Cells.Item(3,7).Formula = '="0.0241"'

$q3.Cells.Item(4,2).Formula = '="001614"'
$q3.Cells.Item(4,3).Formula = '="东方区域发展混合"'
$q3.Cells.Item(4,4).Formula = '="0.22"'
$q3.Cells.Item(4,5).Formula = '="99.06"'
$q3.Cells.Item(4,6).Formula = '="9.19"'
$q3.Cells.Item(4,7).Formula = '="0.0202"'

$q3.Cells.Item(5,2).Formula = '="012244"'
$q3.Cells.Item(5,3).Formula = '="广发金融地产精选股票A"'
$q3.Cells.Item(5,4).Formula = '="0.31"'
$q3.Cells.Item(5,5).Formula = '="85.52"'
$q3.Cells.Item(5,6).Formula = '="4.62"'
$q3.Cells.Item(5,7).Formula = '="0.0143"'

$q3.Cells.Item(6,2).Formula = '="007084"'
$q3.Cells.Item(6,3).Formula = '="天治转型升级混合"'
$q3.Cells.Item(6,4).Formula = '="0.11"'
$q3.Cells.Item(6,5).Formula = '="91.86"'
$q3.Cells.Item(6,6).Formula = '="3.82"'
$q3.Cells.Item(6,7).Formula = '="0.0042"'

$textRange = $q3.Range("B2:G6")
$textRange.Copy()
$textRange.PasteSpecial(-4163)   # xlPasteValues: collapse formula -> literal text
